$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("x86")

# --- Clean up styling on the existing "x86" sheet ---------------------
# The extra number-format cell style (s="1") applied to B:D (and the
# now-unused E column) is no longer needed; strip it back to Normal and
# drop the now-empty E column / empty D cells entirely.
$ws1.Range("B2:D11").Style = "Normal"
[void]$ws1.Range("D4").Clear()
[void]$ws1.Range("D9").Clear()
[void]$ws1.Range("E2:E11").Clear()

# Update the selection on x86 before switching sheets, so the saved
# sheetView reflects the new selected range.
[void]$ws1.Range("D2:D11").Select()

# --- Add the new "arm" worksheet with the ARM NPB results -------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "arm"

$ws2.Range("A1").Value = "class"
$ws2.Range("B1").Value = "Speedup_A"
$ws2.Range("C1").Value = "Speedup_B"
$ws2.Range("D1").Value = "Speedup_C"
$ws2.Range("F1").Value = "Init_A"
$ws2.Range("G1").Value = "Aligned_A"
$ws2.Range("H1").Value = "Init_B"
$ws2.Range("I1").Value = "Aligned_B"
$ws2.Range("J1").Value = "Init_C"
$ws2.Range("K1").Value = "Aligned_C"

$ws2.Range("A2").Value = "bt"
$ws2.Range("B2").Formula = "=F2/G2"
$ws2.Range("C2").Formula = "=H2/I2"
$ws2.Range("D2").Formula = "=J2/K2"
$ws2.Range("F2").Value = 260.43599999999998
$ws2.Range("G2").Value = 264.05399999999997
$ws2.Range("H2").Value = 1154.0840000000001
$ws2.Range("I2").Value = 1132.0039999999999
$ws2.Range("J2").Value = 4722.4179999999997
$ws2.Range("K2").Value = 4643.634

$ws2.Range("A3").Value = "lu"
$ws2.Range("B3").Formula = "=F3/G3"
$ws2.Range("C3").Formula = "=H3/I3"
$ws2.Range("D3").Formula = "=J3/K3"
$ws2.Range("F3").Value = 200.44399999999999
$ws2.Range("G3").Value = 198.184
$ws2.Range("H3").Value = 857.07600000000002
$ws2.Range("I3").Value = 845.43600000000004
$ws2.Range("J3").Value = 3554.2659999999901
$ws2.Range("K3").Value = 3506.1439999999998

$ws2.Range("A4").Value = "mg"
$ws2.Range("B4").Formula = "=F4/G4"
$ws2.Range("C4").Formula = "=H4/I4"
$ws2.Range("F4").Value = 8.8759999999999994
$ws2.Range("G4").Value = 8.8159999999999901
$ws2.Range("H4").Value = 36.576000000000001
$ws2.Range("I4").Value = 36.775999999999897
$ws2.Range("J4").Value = 296.02999999999997

$ws2.Range("A5").Value = "sp"
$ws2.Range("B5").Formula = "=F5/G5"
$ws2.Range("C5").Formula = "=H5/I5"
$ws2.Range("D5").Formula = "=J5/K5"
$ws2.Range("F5").Value = 203.744
$ws2.Range("G5").Value = 179.66200000000001
$ws2.Range("H5").Value = 895.49
$ws2.Range("I5").Value = 798.147999999999
$ws2.Range("J5").Value = 3637.6439999999998
$ws2.Range("K5").Value = 3293.5680000000002

$ws2.Range("A6").Value = "ua"
$ws2.Range("B6").Formula = "=F6/G6"
$ws2.Range("C6").Formula = "=H6/I6"
$ws2.Range("D6").Formula = "=J6/K6"
$ws2.Range("F6").Value = 372.13200000000001
$ws2.Range("G6").Value = 371.36399999999998
$ws2.Range("H6").Value = 1543.7539999999999
$ws2.Range("I6").Value = 1539.252
$ws2.Range("J6").Value = 6080.4459999999999
$ws2.Range("K6").Value = 6121.3559999999998

$ws2.Range("A7").Value = "cg"
$ws2.Range("B7").Formula = "=F7/G7"
$ws2.Range("C7").Formula = "=H7/I7"
$ws2.Range("D7").Formula = "=J7/K7"
$ws2.Range("F7").Value = 6.1440000000000001
$ws2.Range("G7").Value = 5.7880000000000003
$ws2.Range("H7").Value = 271.678
$ws2.Range("I7").Value = 266.65800000000002
$ws2.Range("J7").Value = 720.31999999999903
$ws2.Range("K7").Value = 717.28799999999899

$ws2.Range("A8").Value = "ep"
$ws2.Range("B8").Formula = "=F8/G8"
$ws2.Range("C8").Formula = "=H8/I8"
$ws2.Range("D8").Formula = "=J8/K8"
$ws2.Range("F8").Value = 54.362000000000002
$ws2.Range("G8").Value = 54.5
$ws2.Range("H8").Value = 216.08599999999899
$ws2.Range("I8").Value = 217.16
$ws2.Range("J8").Value = 863.31600000000003
$ws2.Range("K8").Value = 866.31200000000001

$ws2.Range("A9").Value = "ft"
$ws2.Range("B9").Formula = "=F9/G9"
$ws2.Range("C9").Formula = "=H9/I9"
$ws2.Range("F9").Value = 27.945999999999898
$ws2.Range("G9").Value = 28.404
$ws2.Range("H9").Value = 362.92599999999999
$ws2.Range("I9").Value = 351.452

$ws2.Range("A10").Value = "is"
$ws2.Range("B10").Formula = "=F10/G10"
$ws2.Range("C10").Formula = "=H10/I10"
$ws2.Range("D10").Formula = "=J10/K10"
$ws2.Range("F10").Value = 1.722
$ws2.Range("G10").Value = 1.72
$ws2.Range("H10").Value = 6.56
$ws2.Range("I10").Value = 7.45399999999999
$ws2.Range("J10").Value = 28.504000000000001
$ws2.Range("K10").Value = 30.065999999999999

$ws2.Range("A11").Value = "Geomean"
$ws2.Range("B11").Formula = "=F11/G11"
$ws2.Range("C11").Formula = "=H11/I11"
$ws2.Range("D11").Formula = "=J11/K11"
$ws2.Range("F11").Formula = "=GEOMEAN(F1:F10)"
$ws2.Range("G11").Formula = "=GEOMEAN(G1:G10)"
$ws2.Range("H11").Formula = "=GEOMEAN(H1:H10)"
$ws2.Range("I11").Formula = "=GEOMEAN(I1:I10)"
$ws2.Range("J11").Formula = "=GEOMEAN(J1:J10)"
$ws2.Range("K11").Formula = "=GEOMEAN(K1:K10)"

# arm becomes the active/selected sheet, mirroring the selection made on x86
[void]$ws2.Range("D2:D11").Select()

Write-Output "arm sheet added; x86 restyled"
